# Weekly update: insert the newest onion (Cebolla) price week at rows 388-389,
# pushing the previously-last week (old rows 388-393) down to rows 390-395.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows before row 388; existing rows 388:393 shift to 390:395.
$ws.Rows("388:389").Insert()

# --- New row 388 ---
$ws.Cells.Item(388, 1).Value = 11
$ws.Cells.Item(388, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(388, 3).Value = "Bíobío"
$ws.Cells.Item(388, 4).Value = 44628
$ws.Cells.Item(388, 5).Value = 8
$ws.Cells.Item(388, 6).Value = 100112004
$ws.Cells.Item(388, 7).Value = "Cebolla"
$ws.Cells.Item(388, 8).Value = "Sin especificar"
$ws.Cells.Item(388, 9).Value = "1a (cosecha)"
$ws.Cells.Item(388, 10).Value = 270
$ws.Cells.Item(388, 11).Value = 4500
$ws.Cells.Item(388, 12).Value = 5000
$ws.Cells.Item(388, 13).Value = 4722
$ws.Cells.Item(388, 14).Value = "$/malla 18 kilos"
$ws.Cells.Item(388, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(388, 16).Value = 262
$ws.Cells.Item(388, 17).Value = 18
$ws.Cells.Item(388, 18).Value = "Hortaliza"

# --- New row 389 ---
$ws.Cells.Item(389, 1).Value = 11
$ws.Cells.Item(389, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(389, 3).Value = "Bíobío"
$ws.Cells.Item(389, 4).Value = 44628
$ws.Cells.Item(389, 5).Value = 8
$ws.Cells.Item(389, 6).Value = 100112004
$ws.Cells.Item(389, 7).Value = "Cebolla"
$ws.Cells.Item(389, 8).Value = "Sin especificar"
$ws.Cells.Item(389, 9).Value = "2a (cosecha)"
$ws.Cells.Item(389, 10).Value = 150
$ws.Cells.Item(389, 11).Value = 4000
$ws.Cells.Item(389, 12).Value = 4000
$ws.Cells.Item(389, 13).Value = 4000
$ws.Cells.Item(389, 14).Value = "$/malla 18 kilos"
$ws.Cells.Item(389, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(389, 16).Value = 222
$ws.Cells.Item(389, 17).Value = 18
$ws.Cells.Item(389, 18).Value = "Hortaliza"

# Make sure the date cells keep the date number format used throughout column D.
$ws.Cells.Item(388, 4).NumberFormat = $ws.Cells.Item(390, 4).NumberFormat
$ws.Cells.Item(389, 4).NumberFormat = $ws.Cells.Item(390, 4).NumberFormat
